# Apply the changes described by the diff:
#  - Add two new header columns BC1/BD1 with the same style as the other
#    header cells (Odd_CS_3-3_HT, Odd_CS_4-4_HT)
#  - Replace the row-2 match data with the new match's data (and extend it
#    with the two new BC2/BD2 values)
#  - The used-range dimension naturally grows to A1:BD2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells, copying the style of the existing header row ---
$ws.Range("BB1").Copy() | Out-Null
$ws.Range("BC1:BD1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("BC1").Value = "Odd_CS_3-3_HT"
$ws.Range("BD1").Value = "Odd_CS_4-4_HT"

# --- Row 2: new match values ---
# (B2 "11/11/2024" is unchanged by the diff, so it is intentionally left
#  untouched to avoid Excel's automatic text->date coercion)
$ws.Range("A2").Value  = "8IFiFYon"
$ws.Range("C2").Value  = "13:00"
$ws.Range("D2").Value  = "ROMANIA - LIGA 1"
$ws.Range("E2").Value  = "Gloria Buzau"
$ws.Range("F2").Value  = "Petrolul"
$ws.Range("G2").Value  = 3.25
$ws.Range("H2").Value  = 2.9
$ws.Range("I2").Value  = 2.35
$ws.Range("J2").Value  = 4
$ws.Range("K2").Value  = 1.91
$ws.Range("L2").Value  = 3.2
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("O2").Value  = 1.5
$ws.Range("P2").Value  = 2.5
$ws.Range("Q2").Value  = 2.5
$ws.Range("R2").Value  = 1.5
$ws.Range("S2").Value  = 1.57
$ws.Range("T2").Value  = 2.25
$ws.Range("U2").Value  = 2.1
$ws.Range("V2").Value  = 1.67
$ws.Range("W2").Value  = 7.5
$ws.Range("X2").Value  = 15
$ws.Range("Y2").Value  = 13
$ws.Range("Z2").Value  = 34
$ws.Range("AA2").Value = 34
$ws.Range("AB2").Value = 41
$ws.Range("AC2").Value = 6.5
$ws.Range("AD2").Value = 6
$ws.Range("AE2").Value = 19
$ws.Range("AF2").Value = 67
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 6
$ws.Range("AI2").Value = 10
$ws.Range("AJ2").Value = 10
$ws.Range("AK2").Value = 23
$ws.Range("AL2").Value = 23
$ws.Range("AM2").Value = 41
$ws.Range("AN2").Value = 5
$ws.Range("AO2").Value = 21
$ws.Range("AP2").Value = 34
$ws.Range("AQ2").Value = 67
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 351
$ws.Range("AT2").Value = 2.25
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("AW2").Value = 4.33
$ws.Range("AX2").Value = 15
$ws.Range("AY2").Value = 29
$ws.Range("AZ2").Value = 51
$ws.Range("BA2").Value = 81
$ws.Range("BB2").Value = 301
$ws.Range("BC2").Value = 51
$ws.Range("BD2").Value = 51
